$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original stat values (runs, balls, fours, sixes) for rows 2-4
$oldRow2 = @($ws.Range("C2").Text, $ws.Range("D2").Text, $ws.Range("E2").Text, $ws.Range("F2").Text)
$oldRow3 = @($ws.Range("C3").Text, $ws.Range("D3").Text, $ws.Range("E3").Text, $ws.Range("F3").Text)
$oldRow4 = @($ws.Range("C4").Text, $ws.Range("D4").Text, $ws.Range("E4").Text, $ws.Range("F4").Text)

# Keep values stored as text (matching the source data's text-number style)
$ws.Range("C2:F4").NumberFormat = "@"

# Shift rows up cyclically: new row2 = old row3, new row3 = old row4, new row4 = old row2
$ws.Range("C2").Value = $oldRow3[0]
$ws.Range("D2").Value = $oldRow3[1]
$ws.Range("E2").Value = $oldRow3[2]
$ws.Range("F2").Value = $oldRow3[3]

$ws.Range("C3").Value = $oldRow4[0]
$ws.Range("D3").Value = $oldRow4[1]
$ws.Range("E3").Value = $oldRow4[2]
$ws.Range("F3").Value = $oldRow4[3]

$ws.Range("C4").Value = $oldRow2[0]
$ws.Range("D4").Value = $oldRow2[1]
$ws.Range("E4").Value = $oldRow2[2]
$ws.Range("F4").Value = $oldRow2[3]
